# Update the practice sheet's division problems to the new set of
# values, cell-by-cell, while leaving every other part of the
# document (run formatting, paragraph marks, empty answer rows, etc.)
# untouched.
#
# Addressing cells directly via the Tables/Cell object model (rather
# than a document-wide Find & Replace) avoids ambiguity from values
# that are reused between the "old" and "new" sets (e.g. "28÷5=",
# "88÷9=" and "94÷6=" each appear as both an old value in one cell and
# a new value in another).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "47÷5="
$t.Cell(1, 2).Range.Text = "66÷7="
$t.Cell(1, 3).Range.Text = "33÷9="
$t.Cell(1, 4).Range.Text = "28÷5="
$t.Cell(1, 5).Range.Text = "82÷6="

# Row 5
$t.Cell(5, 1).Range.Text = "11÷6="
$t.Cell(5, 2).Range.Text = "23÷2="
$t.Cell(5, 3).Range.Text = "38÷7="
$t.Cell(5, 4).Range.Text = "40÷9="
$t.Cell(5, 5).Range.Text = "47÷7="

# Row 9
$t.Cell(9, 1).Range.Text = "76÷2="
$t.Cell(9, 2).Range.Text = "78÷6="
$t.Cell(9, 3).Range.Text = "88÷9="
$t.Cell(9, 4).Range.Text = "75÷5="
$t.Cell(9, 5).Range.Text = "42÷8="

# Row 13
$t.Cell(13, 1).Range.Text = "94÷6="
$t.Cell(13, 2).Range.Text = "15÷3="
$t.Cell(13, 3).Range.Text = "56÷9="
$t.Cell(13, 4).Range.Text = "59÷3="
$t.Cell(13, 5).Range.Text = "45÷3="

# Row 17
$t.Cell(17, 1).Range.Text = "36÷3="
$t.Cell(17, 2).Range.Text = "86÷9="
$t.Cell(17, 3).Range.Text = "45÷3="
$t.Cell(17, 4).Range.Text = "66÷7="
$t.Cell(17, 5).Range.Text = "72÷5="
